# This workbook report is rebuilt with an additional leading "quantity" column
# (previously the trailing "quantity" column moved to directly follow "Brands"),
# the brand rows re-sorted alphabetically with several brands newly added
# (BUDDEEZ, POPPY & PINE, POPSOCKET, PURELUXE), updated figures for some existing
# brands, and a new trailing "Total" row. All source data is plain text (as in the
# original workbook), so every data value is written with a leading apostrophe to
# keep Excel from reinterpreting number/currency/percent-looking text as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels (column B is now "quantity"; C..N shift right by one
# relative to the original; O/P stay as Gross Profit / Gross Margin).
$headers = @(
    "Brands",
    "quantity",
    "Amazon Top-line Sales (ATS)",
    "Product Cost",
    "Referral Fee",
    "FBA Fulfillment Fee",
    "other transaction fees",
    "Shipping/Kitting Fees",
    "Adjusted Shipping/Kitting Fees",
    "Labeling/Polybagging Fees",
    "Storage Fees",
    "Allocated fees (Premium Services Fee + Subscription)",
    "Cost of Returns",
    "Adjusted Cost of Returns",
    "Gross Profit",
    "Gross Margin"
)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Rows 2-15: brand name in column A (plain text), quantity/money/percent figures
# in columns B-P (forced to text with a leading apostrophe).
$dataRows = @(
    @("BIHLERFLEX", "1", "`$22.99", "`$-14.11", "`$-3.45", "`$-4.99", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-0.20", "`$0.00", "`$0.00", "`$0.00", "`$0.24", "1.04%"),
    @("BUDDEEZ", "0", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-2.20", "`$0.00", "`$0.00", "`$0.00", "`$-2.20", "0.00%"),
    @("CONCORD", "0", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-1.55", "`$0.00", "`$0.00", "`$0.00", "`$-1.55", "0.00%"),
    @("CRICUT", "190", "`$2,060.64", "`$-441.25", "`$-309.43", "`$-740.32", "`$0.00", "`$-18.91", "`$0.00", "`$-0.55", "`$-173.97", "`$-0.18", "`$-70.71", "`$-17.68", "`$377.26", "18.31%"),
    @("HY-C", "1", "`$94.99", "`$-28.67", "`$-14.25", "`$-23.67", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-597.02", "`$-0.01", "`$-152.32", "`$-76.16", "`$-644.79", "-678.80%"),
    @("LITANIA", "85", "`$22,147.35", "`$-11,708.80", "`$-3,567.35", "`$-529.72", "`$0.00", "`$-277.42", "`$0.00", "`$0.00", "`$-90.33", "`$-2.14", "`$-216.63", "`$-54.15", "`$6,194.86", "27.97%"),
    @("MWI-VETONE", "1", "`$24.99", "`$-13.90", "`$-3.75", "`$-3.68", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$3.66", "14.65%"),
    @("PATRICIA NASH", "46", "`$7,216.00", "`$-3,312.00", "`$-1,062.72", "`$-267.91", "`$0.00", "`$-59.47", "`$0.00", "`$-5.61", "`$-118.39", "`$-0.70", "`$-317.24", "`$-317.24", "`$2,131.43", "29.54%"),
    @("POPPY & PINE", "0", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-0.35", "`$0.00", "`$0.00", "`$0.00", "`$-0.35", "0.00%"),
    @("POPSOCKET", "0", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-0.02", "`$0.00", "`$0.00", "`$0.00", "`$-0.02", "0.00%"),
    @("PURELUXE", "0", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$0.00", "`$-0.17", "`$0.00", "`$0.00", "`$0.00", "`$-0.17", "0.00%"),
    @("VIROX", "7,723", "`$271,751.40", "`$-98,771.67", "`$-39,128.88", "`$-55,460.44", "`$0.00", "`$-13,830.67", "`$-13,830.67", "`$-1,301.85", "`$-1,346.00", "`$-26.31", "`$-3,925.61", "`$-3,925.61", "`$57,959.97", "21.33%"),
    @("WABASH VALLEY FARMS", "2,666", "`$109,799.32", "`$-44,673.16", "`$-16,366.10", "`$-20,540.63", "`$0.00", "`$-2,301.96", "`$-2,301.96", "`$0.00", "`$-3,050.59", "`$-10.59", "`$-1,994.42", "`$-997.20", "`$21,859.09", "19.91%"),
    @("Total", "10,713", "`$413,117.68", "`$-158,963.56", "`$-60,455.93", "`$-77,571.36", "`$0.00", "`$-16,488.43", "`$-16,132.63", "`$-1,308.01", "`$-5,380.79", "`$-39.93", "`$-6,676.93", "`$-5,388.04", "`$87,877.43", "21.27%")
)

$startRow = 2
for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $dataRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    for ($col = 2; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = "'" + $rowValues[$col - 1]
    }
}